# Daily attendance processing - 2025-10-23 15:22:33
# Reorders the "Recorded By" (column G) list so that the literal token
# "System" (exact case) is moved to the front of the comma-separated list,
# keeping the relative order of the remaining tokens unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if (-not ($value -is [string])) { continue }
    if ($value -notmatch 'System') { continue }

    $parts = $value -split ',\s*'
    $idx = [Array]::IndexOf($parts, 'System')

    if ($idx -gt 0) {
        $newParts = @('System')
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $idx) {
                $newParts += $parts[$i]
            }
        }
        $cell.Value2 = [string]::Join(', ', $newParts)
    }
}
